$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("D2").Value = 0.04075
$ws.Range("E2").Value = 0.013275
$ws.Range("F2").Value = 0.094
$ws.Range("I2").Value = 0.0002062020761807941
$ws.Range("J2").Value = 0.0002028196411420839
$ws.Range("K2").Value = 6188.8
$ws.Range("L2").Value = 0.5459180523089137
$ws.Range("M2").Value = 3485
$ws.Range("N2").Value = 0.04331102535900926
$ws.Range("O2").Value = 0.563114012409514
$ws.Range("P2").Value = 3485
$ws.Range("Q2").Value = 0.04331102535900926
$ws.Range("R2").Value = 0.563114012409514
$ws.Range("U2").Value = 45056.1
$ws.Range("V2").Value = 0.5599500400797868
$ws.Range("W2").Value = 0.1213331341460992
$ws.Range("X2").Value = 0.06527297149923544
$ws.Range("Y2").Value = 0.05606016264686371
$ws.Range("Z2").Value = 0.1288075256447035
$ws.Range("AA2").Value = 0
$ws.Range("AB2").Value = 0.04120876098108054
$ws.Range("AC2").Value = -0.04119518474884083
$ws.Range("AD2").Value = 87587.5
$ws.Range("AE2").Value = 22.36195081688214
$ws.Range("AF2").Value = 87609.86195081688
$ws.Range("AG2").Value = 42553.76195081688
$ws.Range("AH2").Value = 0.5212565493864787
$ws.Range("AI2").Value = 0.6295681542395184
$ws.Range("AJ2").Value = 0.3459141860403622
$ws.Range("AK2").Value = 0.4522062000082186
$ws.Range("AN2").Value = 12861.60058737151
$ws.Range("AP2").Value = 6248.716879708792

# --- Row 3 ---
$ws.Range("D3").Value = 0.0727
$ws.Range("E3").Value = 0.0174
$ws.Range("F3").Value = 0.076
$ws.Range("K3").Value = 256.7
$ws.Range("L3").Value = 0.4274771024146544
$ws.Range("M3").Value = 195
$ws.Range("N3").Value = 0.05181071817626273
$ws.Range("O3").Value = 0.7596416049863655
$ws.Range("P3").Value = 195
$ws.Range("Q3").Value = 0.05181071817626273
$ws.Range("R3").Value = 0.7596416049863655
$ws.Range("U3").Value = 1550
$ws.Range("V3").Value = 0.4118287855036268
$ws.Range("W3").Value = 0.1345952181208054
$ws.Range("X3").Value = 0.03897990649721352
$ws.Range("Y3").Value = 0.09561531162359184
$ws.Range("Z3").Value = 0.5160264673025695
$ws.Range("AB3").Value = 0.03633369163183226
$ws.Range("AC3").Value = -0.03633369163183226
$ws.Range("AD3").Value = 699.1
$ws.Range("AF3").Value = 699.1
$ws.Range("AG3").Value = -850.9
$ws.Range("AH3").Value = 0.1566505332974814
$ws.Range("AI3").Value = 0.2364380411255411
$ws.Range("AJ3").Value = -0.2921244163691294
$ws.Range("AK3").Value = -0.6048478817173729

# --- Row 4 ---
$ws.Range("D4").Value = 0.0873
$ws.Range("E4").Value = 0.00915
$ws.Range("F4").ClearContents()
$ws.Range("I4").Value = 0.001523131797617088
$ws.Range("J4").Value = 0.001523131797617088
$ws.Range("K4").Value = 600.7
$ws.Range("L4").Value = 0.4474821215733015
$ws.Range("M4").Value = 437.2
$ws.Range("N4").Value = 0.04684703991427806
$ws.Range("O4").Value = 0.7278175461961045
$ws.Range("P4").Value = 437.2
$ws.Range("Q4").Value = 0.04684703991427806
$ws.Range("R4").Value = 0.7278175461961045
$ws.Range("U4").Value = 3131.3
$ws.Range("V4").Value = 0.3355263862844897
$ws.Range("W4").Value = 0.1628751931889049
$ws.Range("X4").Value = 0.0572177856036345
$ws.Range("Y4").Value = 0.1056574075852704
$ws.Range("Z4").Value = 0.1604268474764674
$ws.Range("AA4").Value = 0.0002443512325828743
$ws.Range("AB4").Value = 0.03980632522108923
$ws.Range("AC4").Value = -0.03956197398850635
$ws.Range("AD4").Value = 9161.200000000001
$ws.Range("AE4").Value = 9.376739374394102
$ws.Range("AF4").Value = 9170.576739374395
$ws.Range("AG4").Value = 6039.276739374395
$ws.Range("AH4").Value = 0.4956244233619527
$ws.Range("AI4").Value = 0.7030830223350774
$ws.Range("AJ4").Value = 0.392880851821439
$ws.Range("AK4").Value = 0.6092847037174539
$ws.Range("AN4").Value = 2337.040816326531
$ws.Range("AP4").Value = 1540.631821268978

# --- Row 5 ---
$ws.Range("D5").Value = 0.0619
$ws.Range("E5").Value = 0.0248
$ws.Range("F5").Value = 0.0398
$ws.Range("K5").Value = 3468.2
$ws.Range("L5").Value = 0.6239228596614316
$ws.Range("M5").Value = 1940.7
$ws.Range("N5").Value = 0.04290062802434285
$ws.Range("O5").Value = 0.5595698056628799
$ws.Range("P5").Value = 1940.7
$ws.Range("Q5").Value = 0.04290062802434285
$ws.Range("R5").Value = 0.5595698056628799
$ws.Range("U5").Value = 29077.5
$ws.Range("V5").Value = 0.6427799306321581
$ws.Range("W5").Value = 0.1376318995519681
$ws.Range("X5").Value = 0.05533543133065909
$ws.Range("Y5").Value = 0.08229646822130905
$ws.Range("Z5").Value = 0.1595337998013971
$ws.Range("AB5").Value = 0.04039071943516348
$ws.Range("AC5").Value = -0.04039071943516348
$ws.Range("AD5").Value = 40731.5
$ws.Range("AF5").Value = 40731.5
$ws.Range("AG5").Value = 11654
$ws.Range("AH5").Value = 0.4737950833211195
$ws.Range("AI5").Value = 0.6123969352714413
$ws.Range("AJ5").Value = 0.2048475069035403
$ws.Range("AK5").Value = 0.3113204270972189

# --- Row 6 ---
$ws.Range("D6").Value = 0.117
$ws.Range("E6").Value = 0.102
$ws.Range("I6").Value = 0.0001562607806178755
$ws.Range("J6").Value = 0.0001557643850954392
$ws.Range("K6").Value = 839.5
$ws.Range("L6").Value = 0.4477810966503094
$ws.Range("M6").Value = 397
$ws.Range("N6").Value = 0.03574773087451376
$ws.Range("O6").Value = 0.4729005360333532
$ws.Range("P6").Value = 397
$ws.Range("Q6").Value = 0.03574773087451376
$ws.Range("R6").Value = 0.4729005360333532
$ws.Range("U6").Value = 3627
$ws.Range("V6").Value = 0.3265919896268549
$ws.Range("W6").Value = 0.1499348109517601
$ws.Range("X6").Value = 0.05239699030470918
$ws.Range("Y6").Value = 0.09753782064705091
$ws.Range("Z6").Value = 0.1743175403208683
$ws.Range("AA6").Value = [double]"2.715246447942947e-05"
$ws.Range("AB6").Value = 0.04053833286176468
$ws.Range("AC6").Value = -0.04051118039728525
$ws.Range("AD6").Value = 8560.6
$ws.Range("AE6").Value = 12.98521144248804
$ws.Range("AF6").Value = 8573.585211442489
$ws.Range("AG6").Value = 4946.585211442489
$ws.Range("AH6").Value = 0.4356676925047367
$ws.Range("AI6").Value = 0.5804224269396011
$ws.Range("AJ6").Value = 0.3081564999584238
$ws.Range("AK6").Value = 0.4438674277973019
$ws.Range("AN6").Value = 2962.145328719723
$ws.Range("AP6").Value = 1711.621180429927

# --- Row 7 ---
$ws.Range("D7").Value = -0.00421
$ws.Range("E7").Value = 0.0207
$ws.Range("F7").Value = 0.112
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 190.2
$ws.Range("L7").Value = 0.6659663865546217
$ws.Range("M7").Value = 71.2
$ws.Range("N7").Value = 0.03917253521126761
$ws.Range("O7").Value = 0.3743427970557309
$ws.Range("P7").Value = 71.2
$ws.Range("Q7").Value = 0.03917253521126761
$ws.Range("R7").Value = 0.3743427970557309
$ws.Range("U7").Value = 2085.6
$ws.Range("V7").Value = 1.147447183098592
$ws.Range("W7").Value = 0.1009018567639257
$ws.Range("X7").Value = 0.08761495651052137
$ws.Range("Y7").Value = 0.01328690025340436
$ws.Range("Z7").Value = 0.06522632805006166
$ws.Range("AA7").Value = 0
$ws.Range("AB7").Value = 0.0418791891003964
$ws.Range("AC7").Value = -0.0418791891003964
$ws.Range("AD7").Value = 4200.2
$ws.Range("AE7").Value = 0
$ws.Range("AF7").Value = 4200.2
$ws.Range("AG7").Value = 2114.6
$ws.Range("AH7").Value = 0.6979627106251455
$ws.Range("AI7").Value = 0.6774406864405412
$ws.Range("AJ7").Value = 0.5377651187630335
$ws.Range("AK7").Value = 0.513938510147041
$ws.Range("AN7").ClearContents()
$ws.Range("AP7").ClearContents()

# --- Row 8 ---
$ws.Range("B8").Value = "Doha Bank Q.P.S.C. (DSM:DHBK)"
$ws.Range("D8").Value = -0.0751
$ws.Range("E8").Value = -0.123
$ws.Range("F8").Value = 0.12
$ws.Range("K8").Value = 194.2
$ws.Range("L8").Value = 0.4365025848505282
$ws.Range("M8").Value = 60.4
$ws.Range("N8").Value = 0.02996180366089587
$ws.Range("O8").Value = 0.3110195674562307
$ws.Range("P8").Value = 60.4
$ws.Range("Q8").Value = 0.02996180366089587
$ws.Range("R8").Value = 0.3110195674562307
$ws.Range("U8").Value = 2910
$ws.Range("V8").Value = 1.443523984324619
$ws.Range("W8").Value = 0.05256177768154383
$ws.Range("X8").Value = 0.1397258369011965
$ws.Range("Y8").Value = -0.0871640592196527
$ws.Range("Z8").Value = 0.04553735926305015
$ws.Range("AB8").Value = 0.04314005542555209
$ws.Range("AC8").Value = -0.04314005542555209
$ws.Range("AD8").Value = 9248.6
$ws.Range("AF8").Value = 9248.6
$ws.Range("AG8").Value = 6338.6
$ws.Range("AH8").Value = 0.8210395490257002
$ws.Range("AI8").Value = 0.709847263796147
$ws.Range("AJ8").Value = 0.7587048895804657
$ws.Range("AK8").Value = 0.6264057713212768

# --- Row 9 ---
$ws.Range("D9").Value = 0.0196
$ws.Range("E9").Value = 0.00307
$ws.Range("F9").Value = 0.118
$ws.Range("K9").Value = 465.9
$ws.Range("L9").Value = 0.4803092783505155
$ws.Range("M9").Value = 288.2
$ws.Range("N9").Value = 0.0589173276637501
$ws.Range("O9").Value = 0.618587679759605
$ws.Range("P9").Value = 288.2
$ws.Range("Q9").Value = 0.0589173276637501
$ws.Range("R9").Value = 0.618587679759605
$ws.Range("U9").Value = 1438.3
$ws.Range("V9").Value = 0.2940346716820672
$ws.Range("W9").Value = 0.0787858290352583
$ws.Range("X9").Value = 0.08669115489444476
$ws.Range("Y9").Value = -0.007905325859186463
$ws.Range("Z9").Value = 0.06331592689295039
$ws.Range("AB9").Value = 0.04502753709432868
$ws.Range("AC9").Value = -0.04502753709432868
$ws.Range("AD9").Value = 11106.3
$ws.Range("AF9").Value = 11106.3
$ws.Range("AG9").Value = 9668
$ws.Range("AH9").Value = 0.6942348683264679
$ws.Range("AI9").Value = 0.6498218391813378
$ws.Range("AJ9").Value = 0.6640292315722959
$ws.Range("AK9").Value = 0.6176451798377308

# --- Row 10 ---
$ws.Range("B10").Value = "Ahli Bank Q.P.S.C. (DSM:ABQK)"
$ws.Range("D10").Value = 0.00306
$ws.Range("E10").Value = -0.00542
$ws.Range("F10").Value = 0.046
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 173.4
$ws.Range("L10").Value = 0.6679506933744221
$ws.Range("M10").Value = 95.3
$ws.Range("N10").Value = 0.04142577700499891
$ws.Range("O10").Value = 0.5495963091118801
$ws.Range("P10").Value = 95.3
$ws.Range("Q10").Value = 0.04142577700499891
$ws.Range("R10").Value = 0.5495963091118801
$ws.Range("U10").Value = 1236.4
$ws.Range("V10").Value = 0.5374483807867855
$ws.Range("W10").Value = 0.108071050171393
$ws.Range("X10").Value = 0.0733281573948364
$ws.Range("Y10").Value = 0.03474289277655655
$ws.Range("Z10").Value = 0.07606880182846426
$ws.Range("AA10").Value = 0
$ws.Range("AB10").Value = 0.04725093908528163
$ws.Range("AC10").Value = -0.04725093908528163
$ws.Range("AD10").Value = 3880
$ws.Range("AE10").Value = 0
$ws.Range("AF10").Value = 3880
$ws.Range("AG10").Value = 2643.6
$ws.Range("AH10").Value = 0.6277809238734731
$ws.Range("AI10").Value = 0.6984447004608295
$ws.Range("AJ10").Value = 0.5346979227766429
$ws.Range("AK10").Value = 0.6121144762434009
$ws.Range("AN10").ClearContents()
$ws.Range("AP10").ClearContents()
